# Rename the three parameter labels (column A, rows 2-4) on both sheets,
# appending the "_0" suffix used by the updated test-input naming scheme.

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("performance_params_0")
$wsScaling = $wb.Worksheets.Item("Scaling")

$renames = @{
    "e_modulus" = "e_modulus_0"
    "tensile_strain_at_break" = "tensile_strain_at_break_0"
    "tensile_yield_strength" = "tensile_yield_strength_0"
}

foreach ($ws in @($wsParams, $wsScaling)) {
    for ($r = 2; $r -le 4; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $old = $cell.Value2
        if ($renames.ContainsKey($old)) {
            $cell.Value = $renames[$old]
        }
    }
}

# Widen column A on the "Scaling" sheet to fit the longer labels
# (target stored width ~24.1640625 characters).
$wsScaling.Columns.Item(1).ColumnWidth = 23.33

# Drop the border that used to box in the renamed parameter labels on the
# "performance_params_0" sheet (A2:A4) - the boxed style is no longer used.
$wsParams.Range("A2:A4").Borders.LineStyle = 0

# The active sheet/tab switches from "Scaling" to "performance_params_0",
# and the remembered selections change on each sheet.
$wsScaling.Range("A2:A4").Select()
$wsParams.Activate()
$wsParams.Range("B14").Select()
